$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The H9 cell was an accidental duplicate of H8's text ("Change Trainer's
# username and ask them to add back their trainees"). Remove it - clearing
# both contents and formatting so the cell drops out of the sheet entirely,
# matching the corrected risk-assessment table.
$ws.Range("H9").Clear()

# Update the active selection left by the editor after making the fix.
$ws.Range("H13").Select()
